$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 with corrected species names
$ws.Range("A3").Value = "Lonchocarpus minimiflorus"
$ws.Range("B3").Value = "Calopogonium galactioides"
$ws.Range("C3").Value = "Calopogonium galactioides"
$ws.Range("D3").Value = "loncmi"
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "loncmi"

# Add row 4
$ws.Range("A4").Value = "Morisonia paradoxa"
$ws.Range("B4").Value = "Combretum spinosum"
$ws.Range("C4").Value = "Combretum spinosum"
$ws.Range("D4").Value = "ste2pa"
$ws.Range("F4").Value = "ste2pa"

# Add row 5
$ws.Range("A5").Value = "Plinia gentryi"
$ws.Range("B5").Value = "Hiraea faginea"
$ws.Range("C5").Value = "Hiraea faginea"
$ws.Range("D5").Value = "plinge"
$ws.Range("F5").Value = "plinge"

# Add row 6
$ws.Range("A6").Value = "Verbesina fuscasiccans"
$ws.Range("B6").Value = "Miconia multiplinervia"
$ws.Range("C6").Value = "Miconia multiplinervia"
$ws.Range("D6").Value = "verbfu"
$ws.Range("F6").Value = "verbfu"
